# Update the build timestamp embedded in the version string from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$cells = @()
$cells += $wsAbout.Range("A2")
$cells += $wsAbout.Range("A6")
$cells += $wsData.Range("S2")
$cells += $wsData.Range("S3")
$cells += $wsData.Range("S4")
$cells += $wsData.Range("S5")
$cells += $wsData.Range("S6")
$cells += $wsData.Range("S7")

foreach ($cell in $cells) {
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains($oldStamp)) {
        $cell.Value = $text.Replace($oldStamp, $newStamp)
    }
}
